$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 80, shifting existing rows 80:185 down to 81:186
$ws.Rows.Item(80).Insert()

# Fill in the new row 80 with data
$ws.Cells.Item(80, 1).Value = 9
$ws.Cells.Item(80, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(80, 3).Value = "Metropolitana"
$ws.Cells.Item(80, 4).Value = 44413
$ws.Cells.Item(80, 5).Value = 13
$ws.Cells.Item(80, 6).Value = 100112044
$ws.Cells.Item(80, 7).Value = "Perejil"
$ws.Cells.Item(80, 8).Value = "Sin especificar"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 124
$ws.Cells.Item(80, 11).Value = 9000
$ws.Cells.Item(80, 12).Value = 10000
$ws.Cells.Item(80, 13).Value = 9500
$ws.Cells.Item(80, 14).Value = "`$/docena de atados"
$ws.Cells.Item(80, 15).Value = "Región Metropolitana"
$ws.Cells.Item(80, 16).Value = 3167
$ws.Cells.Item(80, 17).Value = 3
$ws.Cells.Item(80, 18).Value = "Hortaliza"
